$d = $word.ActiveDocument

# Original paragraph text is "Version 1." (V=0 e=1 r=2 s=3 i=4 o=5 n=6
# <space>=7 1=8 .=9, end=10), followed by the _GoBack bookmark.
# Target paragraph text is "Version 2." but re-run as:
#   "Versi" | "on" | <spellEnd/> | " 2" | <bookmark/> | "."
#
# Work from the end of the paragraph backwards so earlier edits don't
# shift the character offsets used by later ones.

# 1) Add a new run holding the final "." after the bookmark (end of
#    paragraph content, position 10).
$d.Range(10, 10).InsertAfter(".")

# 2) Drop the old trailing "." from " 1." so it becomes " 1".
$d.Range(9, 10).Text = ""

# 3) Bump the version number: "1" -> "2".
$d.Range(8, 9).Text = "2"

# 4) Split the "Version" run into "Versi" + "on" (matching the
#    two-run structure in the target) by forcing a run break via a
#    formatting round-trip on the "Versi" portion.
$d.Range(0, 5).Bold = 1
$d.Range(0, 5).Bold = 0
